# fix parentheses, rerun files/figures
# Update "Percent Matching" values (column A) on sheets ED3A and ED3D
# following a rerun of the underlying analysis.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ED3A")
$ws.Range("A2").Value = 48.82943143812709
$ws.Range("A3").Value = 10.7843137254902
$ws.Range("A4").Value = 10.7843137254902
$ws.Range("A5").Value = 10.45751633986928
$ws.Range("A6").Value = 10.45751633986928
$ws.Range("A7").Value = 11.16838487972508
$ws.Range("A8").Value = 10.13745704467354
$ws.Range("A9").Value = 13.91752577319588
$ws.Range("A10").Value = 12.02749140893471
$ws.Range("A12").Value = 8.196721311475411
$ws.Range("A13").Value = 10.81967213114754
$ws.Range("A14").Value = 11.80327868852459
$ws.Range("A15").Value = 10.81967213114754
$ws.Range("A16").Value = 28.40336134453781
$ws.Range("A17").Value = 8.169934640522875
$ws.Range("A18").Value = 8.496732026143791
$ws.Range("A19").Value = 12.19931271477663
$ws.Range("A20").Value = 12.02749140893471
$ws.Range("A23").Value = 12.45901639344262
$ws.Range("A24").Value = 12.45901639344262
$ws.Range("A25").Value = 8.524590163934425
$ws.Range("A26").Value = 10.16393442622951
$ws.Range("A27").Value = 30.08403361344537
$ws.Range("A28").Value = 31.09243697478992
$ws.Range("A29").Value = 9.508196721311476

$ws = $wb.Worksheets.Item("ED3D")
$ws.Range("A2").Value = 39.46488294314381
$ws.Range("A3").Value = 13.0718954248366
$ws.Range("A4").Value = 12.41830065359477
$ws.Range("A5").Value = 8.496732026143791
$ws.Range("A6").Value = 10.7843137254902
